$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2304712
$ws.Range("J17").Value = 2342491.5
$ws.Range("L17").Value = 7027474.5
$ws.Range("N17").Value = -7027810.5
$ws.Range("H41").Value = 906.375
$ws.Range("I41").Value = 714.2857
$ws.Range("J41").Value = 1055.7778
$ws.Range("K41").Value = 714.2857
$ws.Range("L41").Value = 1055.7778
$ws.Range("M41").Value = -274.2857
$ws.Range("N41").Value = -1935.7778
$ws.Range("H135").Value = 1990.2188
$ws.Range("I135").Value = 2255.8333
$ws.Range("J135").Value = 1193.375
$ws.Range("K135").Value = 20302.4997
$ws.Range("L135").Value = 10740.375
$ws.Range("M135").Value = -17767.4997
$ws.Range("N135").Value = -15810.375
$ws.Range("H137").Value = 1378.1
$ws.Range("I137").Value = 1168.2084
$ws.Range("J137").Value = 2217.6667
$ws.Range("K137").Value = 3504.6252
$ws.Range("L137").Value = 6653.000100000001
$ws.Range("M137").Value = -954.6251999999999
$ws.Range("N137").Value = -11753.0001
$ws.Range("H140").Value = 59300
$ws.Range("J140").Value = 59300
$ws.Range("L140").Value = 59300
$ws.Range("N140").Value = -69660
$ws.Range("H141").Value = 2030.6552
$ws.Range("I141").Value = 1732.6666
$ws.Range("K141").Value = 5197.9998
$ws.Range("M141").Value = -17.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4434.826
$ws.Range("I32").Value = 3334.0508
$ws.Range("J32").Value = 10929.4
$ws.Range("K32").Value = 3334.0508
$ws.Range("L32").Value = 10929.4
$ws.Range("M32").Value = -3047.0508
$ws.Range("N32").Value = -11503.4
$ws.Range("H61").Value = 4439.3438
$ws.Range("I61").Value = 4752.3105
$ws.Range("J61").Value = 1414
$ws.Range("K61").Value = 4752.3105
$ws.Range("L61").Value = 1414
$ws.Range("M61").Value = -4540.3105
$ws.Range("N61").Value = -1838
$ws.Range("H74").Value = 1395.6666
$ws.Range("I74").Value = 1285.409
$ws.Range("J74").Value = 1698.875
$ws.Range("K74").Value = 1285.409
$ws.Range("L74").Value = 1698.875
$ws.Range("M74").Value = -411.4090000000001
$ws.Range("N74").Value = -3446.875
$ws.Range("H77").Value = 1395.6666
$ws.Range("I77").Value = 1285.409
$ws.Range("J77").Value = 1698.875
$ws.Range("K77").Value = 6427.045
$ws.Range("L77").Value = 8494.375
$ws.Range("M77").Value = -2059.045
$ws.Range("N77").Value = -17230.375
$ws.Range("H110").Value = 918.2727
$ws.Range("I110").Value = 843
$ws.Range("J110").Value = 1050
$ws.Range("K110").Value = 843
$ws.Range("L110").Value = 1050
$ws.Range("M110").Value = 1202
$ws.Range("N110").Value = -5140
$ws.Range("H122").Value = 1511908.9
$ws.Range("I122").Value = 1712363.4
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 5137090.199999999
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -5134640.199999999
$ws.Range("N122").Value = -30400
$ws.Range("H123").Value = 49426
$ws.Range("J123").Value = 49426
$ws.Range("L123").Value = 49426
$ws.Range("N123").Value = -59226
$ws.Range("H136").Value = 4439.3438
$ws.Range("I136").Value = 4752.3105
$ws.Range("J136").Value = 1414
$ws.Range("K136").Value = 14256.9315
$ws.Range("L136").Value = 4242
$ws.Range("M136").Value = -11706.9315
$ws.Range("N136").Value = -9342

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7753553.5
$ws.Range("I86").Value = 12347179
$ws.Range("J86").Value = 1810.875
$ws.Range("K86").Value = 12347179
$ws.Range("L86").Value = 1810.875
$ws.Range("M86").Value = -12346056
$ws.Range("N86").Value = -4056.875
$ws.Range("H89").Value = 7753553.5
$ws.Range("I89").Value = 12347179
$ws.Range("J89").Value = 1810.875
$ws.Range("K89").Value = 61735895
$ws.Range("L89").Value = 9054.375
$ws.Range("M89").Value = -61730279
$ws.Range("N89").Value = -20286.375
$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2407.758
$ws.Range("I31").Value = 1203.4286
$ws.Range("J31").Value = 3968.926
$ws.Range("K31").Value = 1203.4286
$ws.Range("L31").Value = 3968.926
$ws.Range("M31").Value = -908.4286
$ws.Range("N31").Value = -4558.925999999999
$ws.Range("H34").Value = 2407.758
$ws.Range("I34").Value = 1203.4286
$ws.Range("J34").Value = 3968.926
$ws.Range("K34").Value = 1203.4286
$ws.Range("L34").Value = 3968.926
$ws.Range("M34").Value = -1001.4286
$ws.Range("N34").Value = -4372.925999999999
$ws.Range("H58").Value = 1123.8889
$ws.Range("I58").Value = 852.89655
$ws.Range("J58").Value = 1438.24
$ws.Range("K58").Value = 852.89655
$ws.Range("L58").Value = 1438.24
$ws.Range("M58").Value = -649.89655
$ws.Range("N58").Value = -1844.24
$ws.Range("H86").Value = 3200
$ws.Range("I86").Value = 3600
$ws.Range("J86").Value = 2400
$ws.Range("K86").Value = 3600
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -2477
$ws.Range("N86").Value = -4646
$ws.Range("H89").Value = 3200
$ws.Range("I89").Value = 3600
$ws.Range("J89").Value = 2400
$ws.Range("K89").Value = 18000
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = -12384
$ws.Range("N89").Value = -23232
$ws.Range("H122").Value = 930.3125
$ws.Range("I122").Value = 632.8570999999999
$ws.Range("J122").Value = 1161.6666
$ws.Range("K122").Value = 1898.5713
$ws.Range("L122").Value = 3484.9998
$ws.Range("M122").Value = 551.4287000000002
$ws.Range("N122").Value = -8384.9998
$ws.Range("H134").Value = 2178.7568
$ws.Range("I134").Value = 2386.6206
$ws.Range("J134").Value = 1425.25
$ws.Range("K134").Value = 7159.861800000001
$ws.Range("L134").Value = 4275.75
$ws.Range("M134").Value = -4624.861800000001
$ws.Range("N134").Value = -9345.75
$ws.Range("H136").Value = 1123.8889
$ws.Range("I136").Value = 852.89655
$ws.Range("J136").Value = 1438.24
$ws.Range("K136").Value = 2558.68965
$ws.Range("L136").Value = 4314.72
$ws.Range("M136").Value = -8.689650000000256
$ws.Range("N136").Value = -9414.720000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 3446.6667
$ws.Range("J106").Value = 3446.6667
$ws.Range("L106").Value = 10340.0001
$ws.Range("N106").Value = -12232.0001
$ws.Range("H107").Value = 517.5
$ws.Range("J107").Value = 573.3333
$ws.Range("L107").Value = 1719.9999
$ws.Range("N107").Value = -5559.9999
$ws.Range("H131").Value = 884.25
$ws.Range("I131").Value = 647.375
$ws.Range("J131").Value = 929.369
$ws.Range("K131").Value = 1942.125
$ws.Range("L131").Value = 2788.107
$ws.Range("M131").Value = 3097.875
$ws.Range("N131").Value = -12868.107

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9681.5
$ws.Range("J57").Value = 9681.5
$ws.Range("L57").Value = 9681.5
$ws.Range("N57").Value = -11321.5
$ws.Range("H80").Value = 2477.7778
$ws.Range("I80").Value = 2287.5
$ws.Range("K80").Value = 2287.5
$ws.Range("M80").Value = -1289.5
$ws.Range("H83").Value = 2477.7778
$ws.Range("I83").Value = 2287.5
$ws.Range("K83").Value = 11437.5
$ws.Range("M83").Value = -6445.5
$ws.Range("H122").Value = 26622138
$ws.Range("I122").Value = 38030710
$ws.Range("J122").Value = 2134.6667
$ws.Range("K122").Value = 114092130
$ws.Range("L122").Value = 6404.000100000001
$ws.Range("M122").Value = -114089680
$ws.Range("N122").Value = -11304.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 570.619
$ws.Range("I16").Value = 570.64703
$ws.Range("J16").Value = 570.5
$ws.Range("K16").Value = 570.64703
$ws.Range("L16").Value = 570.5
$ws.Range("M16").Value = -400.64703
$ws.Range("N16").Value = -910.5
$ws.Range("H122").Value = 3018421
$ws.Range("I122").Value = 5497920.5
$ws.Range("J122").Value = 716028.5600000001
$ws.Range("K122").Value = 16493761.5
$ws.Range("L122").Value = 2148085.68
$ws.Range("M122").Value = -16491311.5
$ws.Range("N122").Value = -2152985.68
$ws.Range("H125").Value = 44000
$ws.Range("J125").Value = 44000
$ws.Range("L125").Value = 44000
$ws.Range("N125").Value = -53840
$ws.Range("H133").Value = 40900
$ws.Range("J133").Value = 40900
$ws.Range("L133").Value = 40900
$ws.Range("N133").Value = -45960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 783.0323
$ws.Range("I113").Value = 784.0625
$ws.Range("J113").Value = 781.93335
$ws.Range("K113").Value = 2352.1875
$ws.Range("L113").Value = 2345.80005
$ws.Range("M113").Value = -182.1875
$ws.Range("N113").Value = -6685.80005
